# "Debug Sessions Done und Doku Update"
# Append the 12.02.2023 / 15.02.2023 journal entries to the end of the
# project documentation, after the last (empty) "Listenabsatz" paragraph.

$d = $word.ActiveDocument

# --- anchor: the last paragraph in the document (currently an empty,
# "Listenabsatz"-styled paragraph right before the section break) -------
$anchor = $d.Paragraphs.Last

# 1) New paragraph: "12.02.2023" -- keeps the Listenabsatz style but
#    removes the list's hanging indent (ind left=0), no bullet.
$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pDate1 = $d.Paragraphs.Last
$pDate1.Range.InsertAfter("12.02.2023")
$pDate1.Range.ParagraphFormat.LeftIndent = 0

# 2) New bulleted paragraph: "Eric präsentiert Video, Teaserfolie, und Poster"
$r = $pDate1.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pItem1 = $d.Paragraphs.Last
$pItem1.Range.InsertAfter("Eric präsentiert Video, Teaserfolie, und Poster")
$pItem1.Range.ListFormat.ApplyBulletDefault()

# 3) New bulleted paragraph (same list): "Abklärung Postersession"
$r = $pItem1.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pItem2 = $d.Paragraphs.Last
$pItem2.Range.InsertAfter("Abklärung Postersession")
$pItem2.Range.ListFormat.ApplyBulletDefault()

# 4) New plain paragraph (no list style): "15.02.2023"
$r = $pItem2.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pDate2 = $d.Paragraphs.Last
$pDate2.Range.InsertAfter("15.02.2023")
$pDate2.Range.ListFormat.RemoveNumbers()
$pDate2.Range.Style = $d.Styles.Item("Standard")

# 5) New bulleted paragraph (new list): "Postersession"
$r = $pDate2.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pItem3 = $d.Paragraphs.Last
$pItem3.Range.InsertAfter("Postersession")
$pItem3.Range.ListFormat.ApplyBulletDefault()

# 6) New bulleted paragraph (same list): "Hinweise für KI-Optimierung bekommen"
$r = $pItem3.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pItem4 = $d.Paragraphs.Last
$pItem4.Range.InsertAfter("Hinweise für KI-Optimierung bekommen")
$pItem4.Range.ListFormat.ApplyBulletDefault()

# 7) New bulleted paragraph (same list): "Debugsessions"
$r = $pItem4.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pItem5 = $d.Paragraphs.Last
$pItem5.Range.InsertAfter("Debugsessions")
$pItem5.Range.ListFormat.ApplyBulletDefault()

# 8) Final trailing empty paragraph (no list style) closing out the doc.
$r = $pItem5.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pEnd = $d.Paragraphs.Last
$pEnd.Range.ListFormat.RemoveNumbers()
$pEnd.Range.Style = $d.Styles.Item("Standard")

"ok"
